$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." paragraph, which sits between a blank
# paragraph and the copyright/footer paragraph that also needs to go.
$rng = $d.Content
$find = $rng.Find
$find.ClearFormatting()
$found = $find.Execute(
    "Ver no Jupiter Salvar em pdf Salvar em docx",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $targetStart = $rng.Start

    # Map the found range back to its Paragraphs collection index so we can
    # reach the paragraph immediately before (blank line) and after
    # (copyright/footer line) it without relying on hard-coded offsets.
    $count = $d.Paragraphs.Count
    $idx = -1
    for ($i = 1; $i -le $count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Start -eq $targetStart) {
            $idx = $i
            break
        }
    }

    if ($idx -gt 1 -and $idx -lt $count) {
        $startPara = $d.Paragraphs.Item($idx - 1)
        $endPara = $d.Paragraphs.Item($idx + 1)
        $delRng = $d.Range($startPara.Range.Start, $endPara.Range.End)
        $delRng.Delete()
    }
}
